$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values for rows 2-29, replacing the old Strike# values.
$kValues = @{
    2  = 2
    3  = 4
    4  = 6
    5  = 6
    6  = 3
    7  = 3
    8  = 5
    9  = 6
    10 = 4
    11 = 4
    12 = 4
    13 = 3
    14 = 3
    15 = 7
    16 = 7
    17 = 5
    18 = 10
    19 = 3
    20 = 5
    21 = 6
    22 = 4
    23 = 3
    24 = 4
    25 = 2
    26 = 7
    27 = 1
    28 = 4
    29 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}

$wb.Save()
